# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with the latest scraped snapshot. Both columns are stored as
# plain text in this sheet (not numbers), so for D-column values that
# look numeric (e.g. "212.58", "4.30", "0.520") we prefix the literal
# with a leading apostrophe - Excel's standard "force text / keep
# exact digits" quote-prefix - then reset .Style back to "Normal" so
# no stray text-number-format style is left attached to the cell.
# Values such as "26.324.15" already fail to parse as a number (more
# than one '.'), so they need no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.324.15'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '1.622.52'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''212.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '''18.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.23%  '
$ws.Range("D11").Value = '''0.0815'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.848.99'
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").Value = '1.625.98'
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").Value = '''4.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '''0.520'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '26.330.57'
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").Value = '''62.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.32%  '
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '''202.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").Value = '''4.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.92%  '
$ws.Range("D22").Value = '''9.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '''6.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("E24").Value = '  +7.70%  '
$ws.Range("D25").Value = '''143.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '''15.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").Value = '''6.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").Value = '''0.0527'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.85%  '
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("D34").Value = '''2.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.02%  '
$ws.Range("D35").Value = '''1.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("D36").Value = '1.179.63'
$ws.Range("E36").Value = '  +4.99%  '
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").Value = '''0.810'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '''0.497'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.63%  '
$ws.Range("D42").Value = '''0.793'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("D43").Value = '''5.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("D44").Value = '1.760.33'
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").Value = '''93.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("E46").Value = '  +15.46%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").Value = '''54.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("E51").Value = '  -0.46%  '
